$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Profesor Asociado): "when" changes from "2017 - Actualmente" to 2019
$ws.Range("B2").Value = 2019

# Update row 3 (Profesor Asistente): "when" changes from 2017 to "2017-2018",
# and "why" changes to "Metodos cuantitativos II"
$ws.Range("B3").Value = "2017-2018"
$ws.Range("E3").Value = "Métodos cuantitativos II (Maestría en Psicología)."

# Insert a new row after row 3 for the 2017 "Metodos cuantitativos I" entry
$ws.Rows("4:4").Insert()

$ws.Range("B4").Value = 2017
$ws.Range("E4").Value = "Métodos cuantitativos I (Maestría en Psicología)."

# selection, matching final workbook view state
$ws.Range("C18").Select()
